$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells
$ws.Range("E1").Value = "Execution Time (ms)"
$ws.Range("F1").Value = "Memory Usage (B)"

# Copy the header formatting from an existing header cell (A1) onto the
# two new header cells so they match the rest of the header row.
$ws.Range("A1").Copy()
$ws.Range("E1:F1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# New data values (execution time / memory usage per model)
$ws.Range("E2").Value = 6.878599990159273
$ws.Range("F2").Value = 0

$ws.Range("E3").Value = 7.01979998848401
$ws.Range("F3").Value = 0

$ws.Range("E4").Value = 45.48119998071343
$ws.Range("F4").Value = 0

$ws.Range("E5").Value = 3.086799988523126
$ws.Range("F5").Value = 0
